$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.674.29"
$ws.Range("E2").Value = "  -1.81%  "
$ws.Range("D3").Value = "'1.741.81"
$ws.Range("E3").Value = "  -2.97%  "
$ws.Range("E4").Value = "  +1.18%  "
$ws.Range("D5").Value = "'330.28"
$ws.Range("E5").Value = "  -1.43%  "
$ws.Range("E6").Value = "  +0.85%  "
$ws.Range("D7").Value = "'0.3881"
$ws.Range("E7").Value = "  +1.29%  "
$ws.Range("D8").Value = "'0.3347"
$ws.Range("E8").Value = "  -2.66%  "
$ws.Range("D9").Value = "'45.37"
$ws.Range("E9").Value = "  -5.73%  "
$ws.Range("E10").Value = "  -4.98%  "
$ws.Range("D11").Value = "'0.07159"
$ws.Range("E11").Value = "  -4.36%  "
$ws.Range("E12").Value = "  +1.24%  "
$ws.Range("D13").Value = "'22.10"
$ws.Range("E13").Value = "  -4.39%  "
$ws.Range("D14").Value = "'6.097"
$ws.Range("E14").Value = "  -5.20%  "
$ws.Range("D15").Value = "'1.740.47"
$ws.Range("E15").Value = "  -2.41%  "
$ws.Range("D16").Value = "'6.937"
$ws.Range("E16").Value = "  -3.22%  "
$ws.Range("D17").Value = "'0.00001046"
$ws.Range("E17").Value = "  -4.09%  "
$ws.Range("D18").Value = "'0.06588"
$ws.Range("E18").Value = "  -0.84%  "
$ws.Range("E19").Value = "  +0.89%  "
$ws.Range("D20").Value = "'78.34"
$ws.Range("E20").Value = "  -5.83%  "
$ws.Range("E21").Value = "  -5.94%  "
$ws.Range("D22").Value = "'6.147"
$ws.Range("E22").Value = "  -5.02%  "
$ws.Range("D23").Value = "'27.688.42"
$ws.Range("E23").Value = "  -1.67%  "
$ws.Range("D24").Value = "'11.47"
$ws.Range("E24").Value = "  -5.72%  "
$ws.Range("D25").Value = "'2.395"
$ws.Range("E25").Value = "  +0.47%  "
$ws.Range("D26").Value = "'154.56"
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  -6.73%  "
$ws.Range("D28").Value = "'2.255"
$ws.Range("E28").Value = "  -8.61%  "
$ws.Range("D29").Value = "'1.938.92"
$ws.Range("E29").Value = "  -2.39%  "
$ws.Range("D30").Value = "'1.263"
$ws.Range("E30").Value = "  -13.71%  "
$ws.Range("D31").Value = "'127.94"
$ws.Range("E31").Value = "  -5.14%  "
$ws.Range("D32").Value = "'4.026"
$ws.Range("E32").Value = "  +1.74%  "
$ws.Range("D33").Value = "'5.747"
$ws.Range("E33").Value = "  -7.66%  "
$ws.Range("D34").Value = "'0.08690"
$ws.Range("E34").Value = "  -1.58%  "
$ws.Range("D35").Value = "'11.91"
$ws.Range("E35").Value = "  -7.61%  "
$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").Value = "'5.078"
$ws.Range("E36").Value = "  -5.55%  "
$ws.Range("B37").Value = "WEMIXTOKEN"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").Value = "'1.504"
$ws.Range("E37").Value = "  -0.77%  "
$ws.Range("B38").Value = "TheSandbox"
$ws.Range("C38").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D38").Value = "'0.6425"
$ws.Range("E38").Value = "  -7.02%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.02244"
$ws.Range("E39").Value = "  -8.15%  "
$ws.Range("D40").Value = "'0.06023"
$ws.Range("E40").Value = "  -5.83%  "
$ws.Range("D41").Value = "'0.2083"
$ws.Range("E41").Value = "  -5.70%  "
$ws.Range("D42").Value = "'1.187"
$ws.Range("E42").Value = "  -4.81%  "
$ws.Range("E43").Value = "  +0.88%  "
$ws.Range("D44").Value = "'7.902"
$ws.Range("E44").Value = "  -6.11%  "
$ws.Range("D45").Value = "'13.52"
$ws.Range("E45").Value = "  -5.70%  "
$ws.Range("E46").Value = "  -1.36%  "
$ws.Range("D47").Value = "'0.5940"
$ws.Range("E47").Value = "  -6.73%  "
$ws.Range("D48").Value = "'125.26"
$ws.Range("E48").Value = "  -5.70%  "
$ws.Range("D49").Value = "'1.965"
$ws.Range("E49").Value = "  -7.08%  "
$ws.Range("D50").Value = "'1.145"
$ws.Range("E50").Value = "  -1.03%  "
$ws.Range("E51").Value = "  -6.67%  "
